$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its literal text representation
# (values such as "26.286.62" or "1.007" must remain text, not be
# reinterpreted as numbers/dates by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.286.62"
$ws.Range("E2").Value = "  -0.16%  "

# Row 3
$ws.Range("D3").Value = "1.691.16"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "217.82"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6
$ws.Range("D6").Value = "0.5359"
$ws.Range("E6").Value = "  +1.49%  "

# Row 7
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "0.2723"
$ws.Range("E8").Value = "  +1.03%  "

# Row 9
$ws.Range("D9").Value = "0.06430"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10
$ws.Range("D10").Value = "21.75"
$ws.Range("E10").Value = "  -1.05%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07695"
$ws.Range("E11").Value = "  +2.32%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.697.45"
$ws.Range("E12").Value = "  +0.87%  "

# Row 13
$ws.Range("D13").Value = "4.521"
$ws.Range("E13").Value = "  +0.11%  "

# Row 14
$ws.Range("D14").Value = "0.5796"
$ws.Range("E14").Value = "  +0.11%  "

# Row 15
$ws.Range("D15").Value = "0.000008375"
$ws.Range("E15").Value = "  -1.60%  "

# Row 16
$ws.Range("D16").Value = "66.92"
$ws.Range("E16").Value = "  +3.21%  "

# Row 17
$ws.Range("D17").Value = "26.325.78"
$ws.Range("E17").Value = "  -0.06%  "

# Row 18
$ws.Range("D18").Value = "4.905"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("E19").Value = "  -0.01%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "10.86"
$ws.Range("E20").Value = "  -0.22%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "194.97"
$ws.Range("E21").Value = "  +2.60%  "

# Row 22
$ws.Range("D22").Value = "6.270"
$ws.Range("E22").Value = "  +1.00%  "

# Row 23
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "148.88"
$ws.Range("E24").Value = "  +2.75%  "

# Row 25
$ws.Range("D25").Value = "0.1288"
$ws.Range("E25").Value = "  +2.51%  "

# Row 26
$ws.Range("E26").Value = "  +1.16%  "

# Row 27
$ws.Range("D27").Value = "15.87"
$ws.Range("E27").Value = "  +0.39%  "

# Row 28
$ws.Range("E28").Value = "  +1.80%  "

# Row 29
$ws.Range("D29").Value = "0.06126"
$ws.Range("E29").Value = "  -6.19%  "

# Row 30
$ws.Range("D30").Value = "1.327"
$ws.Range("E30").Value = "  +0.09%  "

# Row 31
$ws.Range("D31").Value = "3.601"
$ws.Range("E31").Value = "  +0.28%  "

# Row 32
$ws.Range("D32").Value = "3.582"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
$ws.Range("D33").Value = "1.689"
$ws.Range("E33").Value = "  +1.72%  "

# Row 34
$ws.Range("E34").Value = "  +0.35%  "

# Row 35
$ws.Range("D35").Value = "0.6201"
$ws.Range("E35").Value = "  -0.31%  "

# Row 36
$ws.Range("D36").Value = "2.426"
$ws.Range("E36").Value = "  +0.89%  "

# Row 37
$ws.Range("E37").Value = "  +0.83%  "

# Row 38
$ws.Range("D38").Value = "0.01643"
$ws.Range("E38").Value = "  +1.29%  "

# Row 39
$ws.Range("D39").Value = "6.188"
$ws.Range("E39").Value = "  -1.71%  "

# Row 40
$ws.Range("D40").Value = "1.110.41"
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("D41").Value = "0.8795"
$ws.Range("E41").Value = "  +0.51%  "

# Row 43
$ws.Range("D43").Value = "100.95"
$ws.Range("E43").Value = "  +0.39%  "

# Row 44
$ws.Range("D44").Value = "1.842.61"
$ws.Range("E44").Value = "  +0.73%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +3.61%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "57.75"
$ws.Range("E46").Value = "  +1.34%  "

# Row 47
$ws.Range("D47").Value = "1.011"
$ws.Range("E47").Value = "  +0.51%  "

# Row 48
$ws.Range("D48").Value = "8.159"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").Value = "0.05290"
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("D50").Value = "0.4291"
$ws.Range("E50").Value = "  -0.10%  "

# Row 51
$ws.Range("D51").Value = "6.056"
$ws.Range("E51").Value = "  -0.52%  "
